$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# --- Plain numeric value updates (style/type unchanged) ---
# Row 15
$ws.Range("N15").Value = -75

# Row 16
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("L16").Value = -32.558139534883
$ws.Range("M16").Value = -40.816326530612
$ws.Range("N16").Value = -84.491978609625

# Row 17
$ws.Range("I17").Value = 53
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = -20.895522388059
$ws.Range("L17").Value = -8.620689655172
$ws.Range("M17").Value = -1.851851851851
$ws.Range("N17").Value = -68.263473053892

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -37.5
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 19.047619047619
$ws.Range("L18").Value = -9.090909090909
$ws.Range("M18").Value = -18.032786885245
$ws.Range("N18").Value = -79.591836734693

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 16
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 105
$ws.Range("J19").Value = 122
$ws.Range("K19").Value = -13.934426229508
$ws.Range("L19").Value = 15.384615384615
$ws.Range("M19").Value = -11.764705882352
$ws.Range("N19").Value = -13.223140495867

# Row 20
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -10.344827586206
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -23.529411764705
$ws.Range("N20").Value = -88.444444444444

# Row 21
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = -54.545454545454
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = -26.086956521739
$ws.Range("I21").Value = 265
$ws.Range("J21").Value = 304
$ws.Range("K21").Value = -12.828947368421
$ws.Range("L21").Value = -8.304498269896
$ws.Range("M21").Value = -17.1875
$ws.Range("N21").Value = -72.338204592901

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -27.272727272727
$ws.Range("I23").Value = 53
$ws.Range("J23").Value = 63
$ws.Range("K23").Value = -15.873015873015
$ws.Range("L23").Value = -15.873015873015
$ws.Range("M23").Value = 17.777777777777

# Row 24
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -30.769230769230
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 33
$ws.Range("H24").Value = 21.212121212121
$ws.Range("I24").Value = 271
$ws.Range("J24").Value = 294
$ws.Range("K24").Value = -7.823129251700
$ws.Range("L24").Value = -17.629179331307
$ws.Range("M24").Value = 6.692913385826

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 85.714285714285
$ws.Range("I25").Value = 115
$ws.Range("J25").Value = 156
$ws.Range("K25").Value = -26.282051282051
$ws.Range("L25").Value = -35.393258426966

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 83.333333333333
$ws.Range("I26").Value = 103
$ws.Range("J26").Value = 87
$ws.Range("K26").Value = 18.390804597701
$ws.Range("L26").Value = 15.730337078651
$ws.Range("M26").Value = -35.625

# Row 28
$ws.Range("F28").Value = 3

# Row 29
$ws.Range("N29").Value = -90.909090909090

# Row 30
$ws.Range("N30").Value = -92.307692307692

# --- Cells changing between numeric and shared-text ("0" / "***.*") representation ---
# Template cells (row 14 / row 22 are untouched by this edit) used as PasteSpecial(xlPasteFormats) sources
$fmtPasteFormats = -4122

# C15: ('15', None, '1') -> ('13', 's', '20')
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial($fmtPasteFormats)

# C18: ('15', None, '1') -> ('13', 's', '20')
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C18").PasteSpecial($fmtPasteFormats)

# C20: ('13', 's', '20') -> ('15', None, '1')
$ws.Range("C20").Value = 1
$ws.Range("G22").Copy()
$ws.Range("C20").PasteSpecial($fmtPasteFormats)

# C27: ('15', None, '1') -> ('13', 's', '20')
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial($fmtPasteFormats)

# C28: ('15', None, '3') -> ('13', 's', '20')
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial($fmtPasteFormats)

# G33: ('15', None, '1') -> ('13', 's', '20')
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("G33").PasteSpecial($fmtPasteFormats)

# H33: ('14', None, '-100') -> ('13', 's', '21')
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H33").PasteSpecial($fmtPasteFormats)
